$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of "totalizacion" data (Fecha, Servicio, Dieta, Cantidad, Valor Total)
$ws.Range("A2").Value = 45829
$ws.Range("A2").NumberFormat = "yyyy-mm-dd"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD"

$ws.Range("B2").Value = "Cena"
$ws.Range("C2").Value = "Liquida total Miel"
$ws.Range("D2").Value = 26
$ws.Range("E2").Value = 326638
